$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) : update column F ("想去人数") ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6881
$ws1.Range("F3").Value = 93
$ws1.Range("F4").Value = 36
$ws1.Range("F5").Value = 449
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 200
$ws1.Range("F10").Value = 1290
$ws1.Range("F11").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F19").Value = 4986
$ws1.Range("F20").Value = 0
$ws1.Range("F22").Value = 435
$ws1.Range("F24").Value = 190

# --- Sheet "全部类型" (sheet4) : update column F ("想去人数") ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 36
$ws4.Range("F6").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F13").Value = 400
$ws4.Range("F14").Value = 138
$ws4.Range("F15").Value = 0
$ws4.Range("F17").Value = 45
$ws4.Range("F18").Value = 0
$ws4.Range("F20").Value = 4986
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 435
$ws4.Range("F25").Value = 0
